# Apply updated dSF (column F) values to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    8  = -2
    10 = -6
    12 = -5
    14 = 5
    15 = -4
    19 = 5
    20 = 0
    23 = 1
    24 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
